$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (A2: 1.0 -> 2.0)
$ws.Range("A2").Value = 2.0

# Update row 3: A3 (2.0 -> 4.0), C3 (I079693 -> I079692), G3 (33.0 -> 34.0), H3 (2025-04-25 05:40:08 -> 2025-04-27 12:56:56)
$ws.Range("A3").Value = 4.0
$ws.Range("C3").Value = "I079692"
$ws.Range("G3").Value = 34.0
$ws.Range("H3").Value = "2025-04-27 12:56:56"

# Add new row 4
$ws.Range("A4").Value = 3.0
$ws.Range("B4").Value = "NewUser"
$ws.Range("C4").Value = "NEW123"
$ws.Range("D4").Value = 100.0
$ws.Range("E4").Value = 1.0
$ws.Range("F4").Value = 1.0
$ws.Range("G4").Value = 35.0
$ws.Range("H4").Value = "2025-04-27 12:35:31"
